$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 3228.5625
$ws.Range("I53").Value = 847.8570999999999
$ws.Range("J53").Value = 5080.222
$ws.Range("K53").Value = 847.8570999999999
$ws.Range("L53").Value = 5080.222
$ws.Range("M53").Value = -210.8570999999999
$ws.Range("N53").Value = -6354.222
$ws.Range("H112").Value = 1450.3684
$ws.Range("J112").Value = 1475.1177
$ws.Range("L112").Value = 4425.3531
$ws.Range("N112").Value = -6641.3531
$ws.Range("H113").Value = 7112.0835
$ws.Range("I113").Value = 4712.6665
$ws.Range("J113").Value = 7911.8887
$ws.Range("K113").Value = 4712.6665
$ws.Range("L113").Value = 7911.8887
$ws.Range("M113").Value = -1458.6665
$ws.Range("N113").Value = -14419.8887
$ws.Range("H132").Value = 1525
$ws.Range("I132").Value = 1300
$ws.Range("K132").Value = 3900
$ws.Range("M132").Value = -1370
$ws.Range("H137").Value = 2564.818
$ws.Range("I137").Value = 1521.95
$ws.Range("K137").Value = 4565.85
$ws.Range("M137").Value = -2015.85

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3215.3333
$ws.Range("I32").Value = 2374.9077
$ws.Range("K32").Value = 2374.9077
$ws.Range("M32").Value = -2087.9077
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H122").Value = 4252.304
$ws.Range("I122").Value = 3592.4285
$ws.Range("J122").Value = 5278.778
$ws.Range("K122").Value = 10777.2855
$ws.Range("L122").Value = 15836.334
$ws.Range("M122").Value = -8327.2855
$ws.Range("N122").Value = -20736.334
$ws.Range("H132").Value = 3533.4443
$ws.Range("I132").Value = 2686.6667
$ws.Range("K132").Value = 8060.000100000001
$ws.Range("M132").Value = -5530.000100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 17747.5
$ws.Range("J6").Value = 17747.5
$ws.Range("L6").Value = 17747.5
$ws.Range("N6").Value = -17973.5
$ws.Range("H21").Value = 63322.668
$ws.Range("J21").Value = 63322.668
$ws.Range("L21").Value = 63322.668
$ws.Range("N21").Value = -63794.668
$ws.Range("H50").Value = 65000
$ws.Range("J50").Value = 65000
$ws.Range("L50").Value = 65000
$ws.Range("N50").Value = -66148
$ws.Range("H86").Value = 3449.4783
$ws.Range("J86").Value = 5478.8
$ws.Range("L86").Value = 5478.8
$ws.Range("N86").Value = -7724.8
$ws.Range("H89").Value = 3449.4783
$ws.Range("J89").Value = 5478.8
$ws.Range("L89").Value = 27394
$ws.Range("N89").Value = -38626
$ws.Range("H134").Value = 1827.6666
$ws.Range("I134").Value = 1118.0416
$ws.Range("K134").Value = 3354.1248
$ws.Range("M134").Value = -819.1248000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1416.3334
$ws.Range("I16").Value = 669.5714
$ws.Range("K16").Value = 669.5714
$ws.Range("M16").Value = -382.5714
$ws.Range("H31").Value = 40432.355
$ws.Range("I31").Value = 2054.15
$ws.Range("K31").Value = 2054.15
$ws.Range("M31").Value = -1759.15
$ws.Range("H34").Value = 40432.355
$ws.Range("I34").Value = 2054.15
$ws.Range("K34").Value = 2054.15
$ws.Range("M34").Value = -1852.15
$ws.Range("H113").Value = 1416.3334
$ws.Range("I113").Value = 669.5714
$ws.Range("K113").Value = 669.5714
$ws.Range("M113").Value = 1500.4286
$ws.Range("H132").Value = 7881.6665
$ws.Range("I132").Value = 5749.5
$ws.Range("J132").Value = 8947.75
$ws.Range("K132").Value = 17248.5
$ws.Range("L132").Value = 26843.25
$ws.Range("M132").Value = -14718.5
$ws.Range("N132").Value = -31903.25
$ws.Range("H134").Value = 2251.4473
$ws.Range("I134").Value = 1452.9678
$ws.Range("J134").Value = 5787.5713
$ws.Range("K134").Value = 4358.903399999999
$ws.Range("L134").Value = 17362.7139
$ws.Range("M134").Value = -1823.903399999999
$ws.Range("N134").Value = -22432.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 11452432
$ws.Range("I4").Value = 5294369
$ws.Range("K4").Value = 15883107
$ws.Range("M4").Value = -15882995
$ws.Range("H12").Value = 77.25
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H35").Value = 5377
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H113").Value = 1158.3684
$ws.Range("I113").Value = 803.3077
$ws.Range("J113").Value = 1927.6666
$ws.Range("K113").Value = 2409.9231
$ws.Range("L113").Value = 5782.9998
$ws.Range("M113").Value = -239.9231
$ws.Range("N113").Value = -10122.9998
$ws.Range("H141").Value = 7740.36
$ws.Range("I141").Value = 3038.1538
$ws.Range("J141").Value = 12834.417
$ws.Range("K141").Value = 9114.4614
$ws.Range("L141").Value = 38503.251
$ws.Range("M141").Value = -3934.4614
$ws.Range("N141").Value = -48863.251

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 780
$ws.Range("I2").Value = 198.66667
$ws.Range("J2").Value = 1167.5555
$ws.Range("K2").Value = 198.66667
$ws.Range("L2").Value = 1167.5555
$ws.Range("M2").Value = -85.66667000000001
$ws.Range("N2").Value = -1393.5555
$ws.Range("H107").Value = 964.55554
$ws.Range("J107").Value = 765
$ws.Range("L107").Value = 765
$ws.Range("N107").Value = -4605
$ws.Range("H114").Value = 59990.5
$ws.Range("J114").Value = 59990.5
$ws.Range("L114").Value = 59990.5
$ws.Range("N114").Value = -68668.5
$ws.Range("H132").Value = 3939.7646
$ws.Range("I132").Value = 3306.889
$ws.Range("J132").Value = 4651.75
$ws.Range("K132").Value = 9920.667000000001
$ws.Range("L132").Value = 13955.25
$ws.Range("M132").Value = -7390.667000000001
$ws.Range("N132").Value = -19015.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6523.294
$ws.Range("I61").Value = 5336.846
$ws.Range("K61").Value = 5336.846
$ws.Range("M61").Value = -5134.846
$ws.Range("H113").Value = 6523.294
$ws.Range("I113").Value = 5336.846
$ws.Range("K113").Value = 5336.846
$ws.Range("M113").Value = -3166.846
$ws.Range("H132").Value = 5683.1113
$ws.Range("I132").Value = 4493.467
$ws.Range("J132").Value = 11631.333
$ws.Range("K132").Value = 13480.401
$ws.Range("L132").Value = 34893.999
$ws.Range("M132").Value = -10950.401
$ws.Range("N132").Value = -39953.999
$ws.Range("H133").Value = 49900
$ws.Range("J133").Value = 49900
$ws.Range("L133").Value = 49900
$ws.Range("N133").Value = -54960
$ws.Range("H136").Value = 4424.1304
$ws.Range("I136").Value = 2855.5
$ws.Range("J136").Value = 5630.769
$ws.Range("K136").Value = 8566.5
$ws.Range("L136").Value = 16892.307
$ws.Range("M136").Value = -6016.5
$ws.Range("N136").Value = -21992.307

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 26969
$ws.Range("J18").Value = 26969
$ws.Range("L18").Value = 26969
$ws.Range("N18").Value = -27315

Write-Host "Applied all Lamia_Profits updates"
